{"js": "// Remove the \"Ver no Jupiter / Salvar em pdf / Salvar em docx\" line, the\n// \"\u00a9 2020 ... Creative Commons Attribution\" copyright line, and the blank\n// paragraph that separated them from the bibliography text, restoring the\n// document to end with the bibliography paragraph followed by the original\n// trailing blank paragraph / page-break paragraph.\n\nconst body = context.document.body;\n\n// Locate the \"Ver no Jupiter Salvar em pdf Salvar em docx\" paragraph.\nconst jupiterResults = body.search(\"Ver no Jupiter Salvar em pdf Salvar em docx\", { matchCase: false });\njupiterResults.load(\"items\");\nawait context.sync();\n\n// Locate the copyright / footer paragraph via a distinctive substring.\nconst copyrightResults = body.search(\"Contact: luizeleno@usp.br\", { matchCase: false });\ncopyrightResults.load(\"items\");\nawait context.sync();\n\nif (jupiterResults.items.length > 0 && copyrightResults.items.length > 0) {\n  const jupiterPara = jupiterResults.items[0].paragraphs.getFirst();\n  const copyrightPara = copyrightResults.items[0].paragraphs.getFirst();\n  // The blank paragraph that sits right before the \"Ver no Jupiter\" line.\n  const blankPara = jupiterPara.previous();\n\n  // Delete bottom-up so earlier deletions don't disturb later references.\n  copyrightPara.delete();\n  jupiterPara.delete();\n  blankPara.delete();\n\n  await context.sync();\n}\n", "ps1": "# Remove the \"Ver no Jupiter Salvar em pdf Salvar em docx\" line, the\n# \"\u00a9 2020 ... Creative Commons Attribution\" copyright line, and the blank\n# paragraph that separated them from the bibliography text. The document\n# should end with the bibliography paragraph followed by the original\n# trailing blank paragraph and the page-break paragraph.\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph indices for the \"Ver no Jupiter...\" line and the\n# copyright line by scanning paragraph text (robust to any shifting caused\n# by earlier edits, unlike hard-coded indices).\n$jupIndex = -1\n$copyIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t -like \"*Ver no Jupiter*\") { $jupIndex = $i }\n    if ($t -like \"*luizeleno@usp.br*\") { $copyIndex = $i }\n}\n\nif ($jupIndex -gt 0 -and $copyIndex -gt 0) {\n    # The blank paragraph sits immediately before the \"Ver no Jupiter\" line.\n    $blankIndex = $jupIndex - 1\n\n    # Delete highest index first so the lower, not-yet-deleted indices stay valid.\n    $indices = @($copyIndex, $jupIndex, $blankIndex) | Sort-Object -Descending\n    foreach ($idx in $indices) {\n        $d.Paragraphs.Item($idx).Range.Delete()\n    }\n}\n"}
